$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plans")

# Integración activación OSM Micronegocio: switch the modeled plan rows
# (5-9) from "Residencial / Sin_TotalPlay_TV" to "Micronegocio /
# M_Sin_TotalPlay_TV" and update the Megas (D column) figures.
$ws.Range("B5").Value = "Micronegocio"
$ws.Range("C5").Value = "M_Sin_TotalPlay_TV"
$ws.Range("D5").Value = 50

$ws.Range("B6").Value = "Micronegocio"
$ws.Range("C6").Value = "M_Sin_TotalPlay_TV"
$ws.Range("D6").Value = 120

$ws.Range("B7").Value = "Micronegocio"
$ws.Range("C7").Value = "M_Sin_TotalPlay_TV"
$ws.Range("D7").Value = 220

$ws.Range("B8").Value = "Micronegocio"
$ws.Range("C8").Value = "M_Sin_TotalPlay_TV"
$ws.Range("D8").Value = 520

$ws.Range("B9").Value = "Micronegocio"
$ws.Range("C9").Value = "M_Sin_TotalPlay_TV"
$ws.Range("D9").Value = 1000

# Leave the cursor where the author left it when saving.
[void]$ws.Activate()
[void]$ws.Range("C11").Select()
